$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header and study-label column (column B) to reflect the new
# citation-style labels, replacing the old "local"/"Michigan_xxxx" values.
$ws.Range("B1").Value = "study"
$ws.Range("B2:B84").Value = "McCoy et al 2021"
$ws.Range("B85:B162").Value = "Kaitany et al 2001"

# Restore the previously-selected cell/view seen in the updated file.
$ws.Range("E90").Select()
